$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.626.10'
$ws.Range('E2').Value = '  -2.72%  '
$ws.Range('D3').Value = '2.904.11'
$ws.Range('E3').Value = '  -3.78%  '
$ws.Range('E4').Value = '  +0.03%  '
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '585.95'
$c.Style = "Normal"
$ws.Range('E5').Value = '  -1.50%  '
$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '147.45'
$c.Style = "Normal"
$ws.Range('E6').Value = '  -0.95%  '
$ws.Range('E7').Value = '  +0.18%  '
$ws.Range('E8').Value = '  -2.71%  '
$ws.Range('D9').Value = '2.900.84'
$ws.Range('E9').Value = '  -3.86%  '
$c = $ws.Range('D10')
$c.NumberFormat = "@"
$c.Value = '6.72'
$c.Style = "Normal"
$ws.Range('E10').Value = '  +4.96%  '
$ws.Range('E11').Value = '  -4.12%  '
$ws.Range('E12').Value = '  -2.55%  '
$ws.Range('E13').Value = '  -3.75%  '
$c = $ws.Range('D14')
$c.NumberFormat = "@"
$c.Value = '34.04'
$c.Style = "Normal"
$ws.Range('E14').Value = '  -1.39%  '
$ws.Range('E15').Value = '  +0.43%  '
$ws.Range('D16').Value = '3.385.29'
$ws.Range('E16').Value = '  -3.73%  '
$c = $ws.Range('D17')
$c.NumberFormat = "@"
$c.Value = '6.82'
$c.Style = "Normal"
$ws.Range('E17').Value = '  -2.76%  '
$ws.Range('D18').Value = '60.532.35'
$ws.Range('E18').Value = '  -2.72%  '
$ws.Range('D19').Value = '2.903.99'
$ws.Range('E19').Value = '  -3.68%  '
$c = $ws.Range('D20')
$c.NumberFormat = "@"
$c.Value = '428.77'
$c.Style = "Normal"
$ws.Range('E20').Value = '  -4.46%  '
$c = $ws.Range('D21')
$c.NumberFormat = "@"
$c.Value = '13.64'
$c.Style = "Normal"
$ws.Range('E21').Value = '  -4.19%  '
$ws.Range('E22').Value = '  -3.15%  '
$ws.Range('E23').Value = '  -4.16%  '
$c = $ws.Range('D24')
$c.NumberFormat = "@"
$c.Value = '80.64'
$c.Style = "Normal"
$ws.Range('E24').Value = '  -2.12%  '
$ws.Range('E25').Value = '  +1.65%  '
$ws.Range('E26').Value = '  -1.91%  '
$c = $ws.Range('D27')
$c.NumberFormat = "@"
$c.Value = '11.83'
$c.Style = "Normal"
$ws.Range('E27').Value = '  -1.78%  '
$ws.Range('E29').Value = '  -0.03%  '
$c = $ws.Range('D30')
$c.NumberFormat = "@"
$c.Value = '7.22'
$c.Style = "Normal"
$ws.Range('E30').Value = '  +0.18%  '
$ws.Range('B31').Value = 'ImmutableX'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$c = $ws.Range('D31')
$c.NumberFormat = "@"
$c.Value = '2.18'
$c.Style = "Normal"
$ws.Range('E31').Value = '  +2.08%  '
$ws.Range('B32').Value = 'PancakeSwap'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$c = $ws.Range('D32')
$c.NumberFormat = "@"
$c.Value = '2.62'
$c.Style = "Normal"
$ws.Range('E32').Value = '  -3.18%  '
$ws.Range('E33').Value = '  -3.78%  '
$ws.Range('E34').Value = '  -2.93%  '
$ws.Range('D35').Value = '0.0₃0838'
$ws.Range('E35').Value = '  -1.62%  '
$ws.Range('E36').Value = '  -2.28%  '
$c = $ws.Range('D37')
$c.NumberFormat = "@"
$c.Value = '5.66'
$c.Style = "Normal"
$ws.Range('E37').Value = '  -3.27%  '
$ws.Range('B38').Value = 'dogwifhat'
$ws.Range('C38').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$c = $ws.Range('D38')
$c.NumberFormat = "@"
$c.Value = '2.97'
$c.Style = "Normal"
$ws.Range('E38').Value = '  -1.18%  '
$ws.Range('B39').Value = 'Stacks'
$ws.Range('C39').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$c = $ws.Range('D39')
$c.NumberFormat = "@"
$c.Value = '2.03'
$c.Style = "Normal"
$ws.Range('E39').Value = '  -1.76%  '
$c = $ws.Range('D40')
$c.NumberFormat = "@"
$c.Value = '49.30'
$c.Style = "Normal"
$ws.Range('E40').Value = '  -1.71%  '
$ws.Range('E41').Value = '  -1.04%  '
$ws.Range('E42').Value = '  -3.60%  '
$ws.Range('E43').Value = '  +1.24%  '
$c = $ws.Range('D44')
$c.NumberFormat = "@"
$c.Value = '41.57'
$c.Style = "Normal"
$ws.Range('E44').Value = '  +1.45%  '
$ws.Range('E45').Value = '  -1.50%  '
$c = $ws.Range('D46')
$c.NumberFormat = "@"
$c.Value = '369.64'
$c.Style = "Normal"
$ws.Range('E46').Value = '  -6.42%  '
$c = $ws.Range('D47')
$c.NumberFormat = "@"
$c.Value = '133.67'
$c.Style = "Normal"
$ws.Range('E47').Value = '  -0.70%  '
$ws.Range('D48').Value = '2.654.05'
$ws.Range('E48').Value = '  -3.03%  '
$c = $ws.Range('D50')
$c.NumberFormat = "@"
$c.Value = '24.98'
$c.Style = "Normal"
$ws.Range('E50').Value = '  +5.05%  '
$ws.Range('E51').Value = '  -1.37%  '
